# update: 1/12/2025: Update last 30 days report
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Vega Security" / Brendan McMenimen job (row 7, Job ID 755) has
# dropped out of the trailing-30-day window, so remove it entirely.
# This shifts every following row up by one.
$ws.Rows(7).Delete()

# Job ID 196 (Blockaid, Brendan McMenimen) advanced to a 2nd Interview.
$ws.Range("E2").Value = "2nd Interview"
$ws.Range("F2").Value = 45993

# Job ID 541 (Blockaid, Lim Yi Jun) advanced to a 2nd Interview.
$ws.Range("E5").Value = "2nd Interview"
$ws.Range("F5").Value = 45995

# Job ID 777 (Adaptive6, SE Director) candidate rows reordered:
# Sean Valois now listed before Itai Heller.
$ws.Range("D12").Value = "Sean Valois"
$ws.Range("E12").Value = "2nd Interview"
$ws.Range("F12").Value = 45978

$ws.Range("D13").Value = "Itai Heller"
$ws.Range("E13").Value = "4th Interview"
$ws.Range("F13").Value = 45989

# Job ID 824 (Blockaid, Adam Palmer) advanced to a 1st Interview.
$ws.Range("E17").Value = "1st Interview"
$ws.Range("F17").Value = 45992
